# Insert a new row at position 249 (pushes the existing row 249..369 down to 250..370)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(249).Insert()

# The row that used to be 249 is now 250; copy its (still intact) values into
# the freshly-inserted row 249, then overwrite the two cells that actually
# differ for the new record (Fecha / Volumen).
$ws.Cells.Item(249, 1).Value = $ws.Cells.Item(250, 1).Value2   # Mercado ID
$ws.Cells.Item(249, 2).Value = $ws.Cells.Item(250, 2).Value2   # Mercado
$ws.Cells.Item(249, 3).Value = $ws.Cells.Item(250, 3).Value2   # Region
$ws.Cells.Item(249, 4).Value = 44839                            # Fecha (new)
$ws.Cells.Item(249, 5).Value = $ws.Cells.Item(250, 5).Value2   # Codreg
$ws.Cells.Item(249, 6).Value = $ws.Cells.Item(250, 6).Value2   # Categoria ID
$ws.Cells.Item(249, 7).Value = $ws.Cells.Item(250, 7).Value2   # Categoria
$ws.Cells.Item(249, 8).Value = $ws.Cells.Item(250, 8).Value2   # Variedad
$ws.Cells.Item(249, 9).Value = $ws.Cells.Item(250, 9).Value2   # Calidad
$ws.Cells.Item(249, 10).Value = 250                             # Volumen (new)
$ws.Cells.Item(249, 11).Value = $ws.Cells.Item(250, 11).Value2 # Precio minimo
$ws.Cells.Item(249, 12).Value = $ws.Cells.Item(250, 12).Value2 # Precio maximo
$ws.Cells.Item(249, 13).Value = $ws.Cells.Item(250, 13).Value2 # Precio promedio ponderado
$ws.Cells.Item(249, 14).Value = $ws.Cells.Item(250, 14).Value2 # Unidad de comercializacion
$ws.Cells.Item(249, 15).Value = $ws.Cells.Item(250, 15).Value2 # Origen
$ws.Cells.Item(249, 16).Value = $ws.Cells.Item(250, 16).Value2 # Precio $/Kg
$ws.Cells.Item(249, 17).Value = $ws.Cells.Item(250, 17).Value2 # Kg o Unidades
$ws.Cells.Item(249, 18).Value = $ws.Cells.Item(250, 18).Value2 # Clasificacion

# Match the date-style used by the rest of the "Fecha" column.
$ws.Cells.Item(249, 4).NumberFormat = $ws.Cells.Item(250, 4).NumberFormat
